{"js": "const replacements = [\n  [\"2024-04-08 Monday\", \"2024-04-09 Tuesday\"],\n  [\"390\u00d73=1170\", \"149\u00d77=1043\"],\n  [\"780\u00d79=7020\", \"731\u00d74=2924\"],\n  [\"789\u00d76=4734\", \"670\u00d79=6030\"],\n  [\"844\u00d73=2532\", \"659\u00d79=5931\"],\n  [\"398\u00d77=2786\", \"982\u00d75=4910\"],\n  [\"927\u00d74=3708\", \"760\u00d73=2280\"],\n  [\"947\u00d73=2841\", \"840\u00d78=6720\"],\n  [\"821\u00d74=3284\", \"162\u00d73=486\"],\n  [\"302\u00d74=1208\", \"556\u00d76=3336\"],\n  [\"128\u00d74=512\", \"300\u00d76=1800\"],\n  [\"775\u00d78=6200\", \"562\u00d79=5058\"],\n  [\"381\u00d72=762\", \"603\u00d77=4221\"],\n  [\"896\u00d77=6272\", \"829\u00d73=2487\"],\n  [\"268\u00d78=2144\", \"148\u00d78=1184\"],\n  [\"116\u00d78=928\", \"783\u00d78=6264\"],\n  [\"499\u00d79=4491\", \"760\u00d78=6080\"],\n  [\"132\u00d72=264\", \"477\u00d77=3339\"],\n  [\"762\u00d76=4572\", \"455\u00d76=2730\"],\n  [\"764\u00d73=2292\", \"192\u00d77=1344\"],\n  [\"225\u00d75=1125\", \"208\u00d76=1248\"],\n  [\"221\u00d76=1326\", \"471\u00d76=2826\"],\n  [\"574\u00d79=5166\", \"574\u00d78=4592\"],\n  [\"697\u00d72=1394\", \"508\u00d77=3556\"],\n  [\"930\u00d79=8370\", \"880\u00d76=5280\"],\n  [\"760\u00d75=3800\", \"108\u00d76=648\"],\n];\n\nconst body = context.document.body;\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}", "ps1": "$d = $word.ActiveDocument\n\n$replacements = @(\n    ,@(\"2024-04-08 Monday\", \"2024-04-09 Tuesday\")\n    ,@(\"390\u00d73=1170\", \"149\u00d77=1043\")\n    ,@(\"780\u00d79=7020\", \"731\u00d74=2924\")\n    ,@(\"789\u00d76=4734\", \"670\u00d79=6030\")\n    ,@(\"844\u00d73=2532\", \"659\u00d79=5931\")\n    ,@(\"398\u00d77=2786\", \"982\u00d75=4910\")\n    ,@(\"927\u00d74=3708\", \"760\u00d73=2280\")\n    ,@(\"947\u00d73=2841\", \"840\u00d78=6720\")\n    ,@(\"821\u00d74=3284\", \"162\u00d73=486\")\n    ,@(\"302\u00d74=1208\", \"556\u00d76=3336\")\n    ,@(\"128\u00d74=512\", \"300\u00d76=1800\")\n    ,@(\"775\u00d78=6200\", \"562\u00d79=5058\")\n    ,@(\"381\u00d72=762\", \"603\u00d77=4221\")\n    ,@(\"896\u00d77=6272\", \"829\u00d73=2487\")\n    ,@(\"268\u00d78=2144\", \"148\u00d78=1184\")\n    ,@(\"116\u00d78=928\", \"783\u00d78=6264\")\n    ,@(\"499\u00d79=4491\", \"760\u00d78=6080\")\n    ,@(\"132\u00d72=264\", \"477\u00d77=3339\")\n    ,@(\"762\u00d76=4572\", \"455\u00d76=2730\")\n    ,@(\"764\u00d73=2292\", \"192\u00d77=1344\")\n    ,@(\"225\u00d75=1125\", \"208\u00d76=1248\")\n    ,@(\"221\u00d76=1326\", \"471\u00d76=2826\")\n    ,@(\"574\u00d79=5166\", \"574\u00d78=4592\")\n    ,@(\"697\u00d72=1394\", \"508\u00d77=3556\")\n    ,@(\"930\u00d79=8370\", \"880\u00d76=5280\")\n    ,@(\"760\u00d75=3800\", \"108\u00d76=648\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n    $range = $d.Content\n    $range.Find.ClearFormatting()\n    $range.Find.Replacement.ClearFormatting()\n    $range.Find.Text = $oldText\n    $range.Find.Replacement.Text = $newText\n    $range.Find.Forward = $true\n    $range.Find.Wrap = 1\n    $range.Find.Execute([ref]$oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2) | Out-Null\n}"}
